$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FEDBANK")

# Row 7
$ws.Range("F7").Value = 151.15
$ws.Range("G7").Value = 154.5
$ws.Range("H7").Value = 149.35
$ws.Range("I7").Value = 153.95
$ws.Range("J7").Value = 153.2

# Row 9
$ws.Range("G9").Value = 151.95
$ws.Range("H9").Value = 149.45
$ws.Range("I9").Value = 150.2

# Row 10
$ws.Range("G10").Value = 151.65
$ws.Range("H10").Value = 149.35
$ws.Range("I10").Value = 151.1

# Row 11
$ws.Range("G11").Value = 152.55
$ws.Range("H11").Value = 150.85
$ws.Range("I11").Value = 151.75

# Row 12
$ws.Range("G12").Value = 151.85
$ws.Range("H12").Value = 150.9
$ws.Range("I12").Value = 151.65

# Row 13
$ws.Range("G13").Value = 152
$ws.Range("H13").Value = 150.4
$ws.Range("I13").Value = 150.45

# Row 14
$ws.Range("G14").Value = 150.75
$ws.Range("H14").Value = 149.95
$ws.Range("I14").Value = 150.55

# Row 15
$ws.Range("G15").Value = 150.95
$ws.Range("H15").Value = 150.55
$ws.Range("I15").Value = 150.85

# Row 16
$ws.Range("G16").Value = 152
$ws.Range("H16").Value = 150.85
$ws.Range("I16").Value = 151.55

# Row 17
$ws.Range("G17").Value = 151.8
$ws.Range("H17").Value = 151.25
$ws.Range("I17").Value = 151.35

# Row 18
$ws.Range("G18").Value = 151.55
$ws.Range("H18").Value = 150.8
$ws.Range("I18").Value = 151.15

# Row 19
$ws.Range("G19").Value = 151.5
$ws.Range("H19").Value = 150.6
$ws.Range("I19").Value = 151.5

# Row 20
$ws.Range("G20").Value = 154.15
$ws.Range("H20").Value = 151.45
$ws.Range("I20").Value = 153.7

# Row 21
$ws.Range("G21").Value = 154.5
$ws.Range("H21").Value = 153.5
$ws.Range("I21").Value = 154.4
